$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-10-08 12:48:28"

for ($row = 2; $row -le 16; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
